# Update the average_county_temperature column (AA) with new NOAA-based
# temperature values. The new value depends on the facility_id (column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of facility_id -> new average_county_temperature value
$facilityTemps = @{
    1000328 = -3.222222222222223
    1000378 = 0.4166666666666667
    1002874 = 12.51681286549706
    1005657 = 12.51681286549706
    1005659 = 12.51681286549706
    1006580 = 0.4166666666666667
    1012424 = 15.74228395061728
    1013538 = 15.74228395061728
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162, col G = facility_id

for ($r = 2; $r -le $lastRow; $r++) {
    $facilityId = $ws.Cells.Item($r, 7).Value()  # column G = facility_id
    if ($null -ne $facilityId -and $facilityTemps.ContainsKey([int]$facilityId)) {
        $ws.Cells.Item($r, 27).Value2 = $facilityTemps[[int]$facilityId]  # column AA = average_county_temperature
    }
}

$wb.Save()
